$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refactored DICOM parameters and how columns are written for values in
# dictionaries: the StudyDescription (col J) and StudyDate (col K) values
# had been written to the wrong columns, so swap J/K for the header row
# and every data row that carries scan/DICOM metadata.
#
# Use Copy/PasteSpecial (via an out-of-range scratch cell) rather than
# Value/Value2 so the swapped cells keep their original text data type
# instead of being reinterpreted as numbers.
$scratch = $ws.Range("Z100")
$rows = @(1, 12, 14, 17, 21, 24, 26)
foreach ($r in $rows) {
    $jCell = $ws.Cells.Item($r, 10)
    $kCell = $ws.Cells.Item($r, 11)

    $kCell.Copy()
    $scratch.PasteSpecial()

    $jCell.Copy()
    $kCell.PasteSpecial()

    $scratch.Copy()
    $jCell.PasteSpecial()
}
$scratch.Clear()

# Update the active selection on the sheet to K1.
$ws.Range("K1").Select()
